# fix: prevent hidden columns from being labeled upon detecting changes
#
# The "Änderung" (change) marker column L was being stamped with "ÄNDERUNG"
# for every row, including rows whose only "difference" came from a column
# that should have been skipped while diffing. This also left a handful of
# "first row of a segment group" rows (93, 96, 100, 103, 107, 111, 115, 121,
# 124) stuck with the old (pre-regeneration) cell formatting instead of the
# lightly-shaded "group header" look used by every other such row above them.
#
# Fix both symptoms for rows 93-126:
#   1. Re-stamp the formatting of the group-header rows to match the
#      existing correctly-formatted group-header rows (e.g. row 2).
#   2. Clear out the erroneous "ÄNDERUNG" marker in column L for every row
#      in the block (header rows and normal data rows alike), and make sure
#      column L's formatting matches the un-flagged look (e.g. L2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the formatting of the group-header rows that still carry the old
#    row styling, by copying the formats from an already-correct header row.
$headerRows = @(93, 96, 100, 103, 107, 111, 115, 121, 124)
$ws.Range("A2:V2").Copy()
foreach ($r in $headerRows) {
    $target = "A" + $r + ":V" + $r
    $ws.Range($target).PasteSpecial(-4122)
}

# 2. Column L ("Änderung") formatting should match the un-flagged style for
#    the whole affected block, and the erroneous change marker text must be
#    removed.
$ws.Range("L2").Copy()
$ws.Range("L93:L126").PasteSpecial(-4122)
$ws.Range("L93:L126").Value = ""

$excel.CutCopyMode = 0
